$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New layout: A=Metodo, B=Rx, C=Ry, D=CL, E=Entropia, F=SSIM
# Write column A top-to-bottom first, then the rest of row 1, then the
# numeric data, so new shared-string entries land in the same order as
# the target workbook (Metodo, method names, then Rx/Ry/CL/Entropia/SSIM).
$ws.Range("A1").Value = "Metodo"
$ws.Range("A2").Value = "SMARTER"
$ws.Range("A3").Value = "Fuzzy"
$ws.Range("A4").Value = "TOPSIS"
$ws.Range("A5").Value = "GRA"
$ws.Range("A6").Value = "CODAS"
$ws.Range("A7").Value = "MABAC"
$ws.Range("A8").Value = "VIKOR"
$ws.Range("A9").Value = "PROMETHEE II"

$ws.Range("B1").Value = "Rx"
$ws.Range("C1").Value = "Ry"
$ws.Range("D1").Value = "CL"
$ws.Range("E1").Value = "Entropia"
$ws.Range("F1").Value = "SSIM"

# Row 2 - SMARTER
$ws.Range("B2").Value = 10
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 0.97880382474159178
$ws.Range("E2").Value = 7.8877697597406691
$ws.Range("F2").Value = 0.64031770130013688

# Row 3 - Fuzzy
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 0.164546402835752
$ws.Range("E3").Value = 7.5093335344058882
$ws.Range("F3").Value = 0.71941541490743555

# Row 4 - TOPSIS
$ws.Range("B4").Value = 9
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 7.8825307012206904
$ws.Range("F4").Value = 0.642942557121568

# Row 5 - GRA
$ws.Range("B5").Value = 10
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 0.97880382474159178
$ws.Range("E5").Value = 7.8877697597406691
$ws.Range("F5").Value = 0.64031770130013688

# Row 6 - CODAS
$ws.Range("B6").Value = 22
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 0.79222740479055287
$ws.Range("E6").Value = 7.9035305827081093
$ws.Range("F6").Value = 0.61578413807230903

# Row 7 - MABAC
$ws.Range("B7").Value = 22
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 0.79222740479055287
$ws.Range("E7").Value = 7.9035305827081093
$ws.Range("F7").Value = 0.61578413807230903

# Row 8 - VIKOR
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = 0.36898499644653843
$ws.Range("E8").Value = 7.6236212627024926
$ws.Range("F8").Value = 0.71355915992509844

# Row 9 - PROMETHEE II
$ws.Range("B9").Value = 22
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = 0.79222740479055287
$ws.Range("E9").Value = 7.9035305827081093
$ws.Range("F9").Value = 0.61578413807230903

# Column widths
$ws.Columns("A").ColumnWidth = 13.28515625
$ws.Columns("B").ColumnWidth = 3.140625
$ws.Columns("C").ColumnWidth = 3.140625
$ws.Columns("D").ColumnWidth = 12
$ws.Columns("E").ColumnWidth = 12
$ws.Columns("F").ColumnWidth = 12
